# DistributedPubSub.pptx: "Upgraded NServiceBus 2.1 to 2.0" edit.
#
# 1) The slide master and every slide layout carry a cached
#    datetimeFigureOut field ("11/3/2009") in their Date placeholder -
#    refresh it to 10/28/2009 everywhere.
# 2) Slide 1 has two message-arrow/label pairs near the bottom of the
#    diagram that get repositioned/rotated (and one loses its end
#    connection site).
# 3) The entrance-wipe animations that play on those two arrows swap
#    their wipe directions with each other.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "today" date field wherever it is used
#    (ppPlaceholderDate = 16) - the slide master plus all 11 layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = "10/28/2009"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Reposition the "Msg 1" / "Msg 2" connectors + labels on slide 1.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# EMU -> points helper. The engine truncates (floors) when it converts
# the COM "points" float back to EMU on save, so nudge by half an EMU
# first to make sure it lands back on the exact integer EMU we want.
function Pt([double]$emu) {
    return ($emu + 0.5) / 12700
}

# "Straight Arrow Connector 135" (msg 1 arrow): was a connector glued
# by its end to shape 128; it becomes a free-floating rotated arrow.
$arrow1 = $slide.Shapes.Item("Straight Arrow Connector 135")
$arrow1.ConnectorFormat.EndDisconnect()
$arrow1.Rotation = 90
$arrow1.VerticalFlip = $false
$arrow1.Left = Pt 5295900
$arrow1.Top = Pt 5372100
$arrow1.Width = Pt 685800
$arrow1.Height = Pt 304800

# "TextBox 136" is the "Msg 1" label that rides along with it.
$label1 = $slide.Shapes.Item("TextBox 136")
$label1.Left = Pt 5105400
$label1.Top = Pt 5300246

# "Straight Arrow Connector 137" (msg 2 arrow): stays unglued, just
# gets rotated/flipped/resized and moved.
$arrow2 = $slide.Shapes.Item("Straight Arrow Connector 137")
$arrow2.Rotation = 270
$arrow2.HorizontalFlip = $true
$arrow2.Left = Pt 6553200
$arrow2.Top = Pt 5181600
$arrow2.Width = Pt 762000
$arrow2.Height = Pt 762000

# "TextBox 138" is the "Msg 2" label that rides along with it.
$label2 = $slide.Shapes.Item("TextBox 138")
$label2.Left = Pt 6934200
$label2.Top = Pt 5334000

# ---------------------------------------------------------------------
# 3) Swap the entrance-wipe directions on the two arrows' animations.
#    (direction 1 = wipe(up), 4 = wipe(left))
# ---------------------------------------------------------------------
$mainSeq = $slide.TimeLine.MainSequence
for ($i = 1; $i -le $mainSeq.Count; $i++) {
    $effect = $mainSeq.Item($i)
    if ($effect.Shape.Name -eq "Straight Arrow Connector 135") {
        $effect.EffectParameters.Direction = 1
    } elseif ($effect.Shape.Name -eq "Straight Arrow Connector 137") {
        $effect.EffectParameters.Direction = 4
    }
}
